# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) timestamps for the first data row
# on the zh-cn and de-de handback-status sheets, reflecting a fresh
# report generation run.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 15:15:23"
$wsZhCn.Range("H2").Value = "2016-03-22 15:15:49"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 15:15:29"
$wsDeDe.Range("H2").Value = "2016-03-22 15:15:56"
